$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.858.86'
$ws.Range("E2").Value = '  +0.20%  '

$ws.Range("D3").Value = '1.883.60'
$ws.Range("E3").Value = '  -0.18%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.010'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.30%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '336.62'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.65%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.010'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.31%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4685'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.70%  '

$ws.Range("E8").Value = '  +0.96%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.08048'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.21%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '45.59'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -4.19%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.015'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -1.10%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '22.07'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -0.28%  '

$ws.Range("D13").Value = '1.888.31'
$ws.Range("E13").Value = '  +0.43%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.018'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +0.43%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.308'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +2.36%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.012'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +0.22%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '89.38'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +2.36%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.06736'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.05%  '

$ws.Range("E19").Value = '  -0.15%  '

$ws.Range("E20").Value = '  +0.36%  '

$ws.Range("E21").Value = '  +0.28%  '

$ws.Range("D22").Value = '27.870.95'
$ws.Range("E22").Value = '  +0.15%  '

$ws.Range("E23").Value = '  -0.14%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '11.04'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.47%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.318'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.59%  '

$ws.Range("D26").Value = '2.105.29'
$ws.Range("E26").Value = '  +0.02%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '159.52'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.14%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '19.88'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -1.50%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.160'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +2.46%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '5.518'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -1.31%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '122.08'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +0.02%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.9875'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +0.64%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.09509'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +0.22%  '

$ws.Range("E34").Value = '  +0.46%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '5.358'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +0.06%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.355'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -6.65%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.06098'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -1.00%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.02252'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -0.50%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '8.343'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +3.16%  '

$ws.Range("E40").Value = '  -1.40%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.009'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +0.22%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.6017'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +0.05%  '

$ws.Range("E43").Value = '  +0.23%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '10.41'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.84%  '

$ws.Range("B45").Value = 'WEMIXTOKEN'
$ws.Range("C45").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.250'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.77%  '

$ws.Range("B46").Value = 'Decentraland'
$ws.Range("C46").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.5691'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.07%  '

$ws.Range("E47").Value = '  -0.32%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.950'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.32%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.06790'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -1.74%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '112.91'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.90%  '

$ws.Range("B51").Value = 'PancakeSwap'
$ws.Range("C51").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '3.045'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -10.49%  '
